$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "sesi" column header
$ws.Range("F1").Value = "sesi"

# Fill in "sesi" values for existing rows
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 2

# Add new guest row (row 5)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Irham"
$ws.Range("C5").Value = "Bogor"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "082111002299"
$ws.Range("E5").Value = "biasa"
$ws.Range("F5").Value = 2

# Update selection to match target state
$ws.Range("E12").Select()
